$d = $word.ActiveDocument

# 1) Zip-file description paragraph: add the ettersendelse.pdf mention after
#    "(sjølve søknaden)" and fix "an ignoreres" -> "kan ignoreres". Both
#    search strings are unique in the document, so a plain document-wide
#    Find/Replace is safe here.
$d.Content.Find.Execute(
    " (sjølve søknaden) og Soknad-j",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " (sjølve søknaden), ettersendelse.pdf(for ettersendelser) og Soknad-j",
    2)

$d.Content.Find.Execute(
    "brukes for maskinell import og an ignoreres: ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "brukes for maskinell import og kan ignoreres: ",
    2)

# 2) "Hvis kommunen ikke allerede har tilgang til Fiks-plattformen ..." —
#    this exact sentence appears twice in the document; only the second
#    occurrence (the one that precedes the "Informasjon" table further
#    down) gets the new "SvarUt" clause. Locate it through the paragraph
#    collection and restrict the Find/Replace to that paragraph's own
#    [Start,End) span via $d.Range(start, end) (Paragraph.Range.Find does
#    not respect the paragraph boundary in this host, but a fresh
#    Document.Range does), so the other occurrence is left untouched.
$count = $d.Paragraphs.Count
$fiksHits = 0
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs($i)
    $t = $p.Range.Text.Trim()
    if ($t -like "Hvis kommunen ikke allerede har tilgang til Fiks-plattformen (det vil si kommuner som ikke har tatt i bruk MinSide), trenger vi følgende informasjon:*") {
        $fiksHits = $fiksHits + 1
        if ($fiksHits -eq 2) {
            $rStart = $p.Range.Start
            $rEnd = $p.Range.End
            $scoped = $d.Range($rStart, $rEnd)
            $scoped.Find.Execute(
                "Fiks-plattformen ",
                $true, $false, $false, $false, $false, $true, 1, $false,
                "Fiks-plattformen, de som kun bruker SvarUt må fylle ut dette. ",
                2)
        }
    }
}

# 3) "til organisasjonsnummeret til NAV-avdelingen " -> add "kommunal "
#    before "NAV-avdelingen". This phrase (or a close variant) also occurs
#    in two other paragraphs elsewhere in the document, so scope the
#    Find/Replace to the one paragraph matching the full sentence from the
#    diff (the "OBS! Denne personen må ha post/arkiv-rollen ..." one).
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs($i)
    $t = $p.Range.Text
    if ($t -like "*OBS!*post/arkiv-rollen*organisasjonsnummeret til NAV-avdelingen*") {
        $rStart = $p.Range.Start
        $rEnd = $p.Range.End
        $scoped = $d.Range($rStart, $rEnd)
        $scoped.Find.Execute(
            "til organisasjonsnummeret til NAV-avdelingen ",
            $true, $false, $false, $false, $false, $true, 1, $false,
            "til organisasjonsnummeret til kommunal NAV-avdelingen ",
            2)
    }
}
